# Insert a new price record as row 235 in the daily "Pepino ensalada" sheet.
# Excel's row Insert pushes the existing rows 235..287 down to 236..288
# (preserving all their values/formatting), matching the diff exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(235).Insert()

$ws.Range("A235").Value2 = 7
$ws.Range("B235").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C235").Value2 = "Ñuble"
$ws.Range("D235").Value2 = 44943
$ws.Range("E235").Value2 = 16
$ws.Range("F235").Value2 = 100112043
$ws.Range("G235").Value2 = "Pepino ensalada"
$ws.Range("H235").Value2 = "Sin especificar"
$ws.Range("I235").Value2 = "Primera"
$ws.Range("J235").Value2 = 60
$ws.Range("K235").Value2 = 9000
$ws.Range("L235").Value2 = 9000
$ws.Range("M235").Value2 = 9000
$ws.Range("N235").Value2 = "$/caja 80 unidades"
$ws.Range("O235").Value2 = "Región del Maule"
$ws.Range("P235").Value2 = 112
$ws.Range("Q235").Value2 = 80
$ws.Range("R235").Value2 = "Hortaliza"

# Keep the date column formatted the same way as its neighbours (style index 2).
$ws.Range("D235").NumberFormat = $ws.Range("D236").NumberFormat
